$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange
$para1 = $tr.Paragraphs(1, 1)

$para1.Text = "Car Sales "
$para1.InsertAfter("dataset ") | Out-Null
$para1.InsertAfter("EDA Analysis and Observations") | Out-Null
